# [FIX] Inventory 버그 수정
# Renames the header row (column-name) cells across all three sheets so
# they no longer carry the leading underscore used as an internal/private
# field marker (e.g. "_id" -> "id", "_itemName" -> "itemName", ...), and
# restores the previously-saved sheet/selection UI state.

$wb = $excel.ActiveWorkbook

# --- ItemDatas ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ItemDatas")
$ws1.Range("A1").Value = "id"
$ws1.Range("B1").Value = "itemName"
$ws1.Range("C1").Value = "itemExplanation"
$ws1.Range("D1").Value = "price"
$ws1.Range("E1").Value = "dropPrefabPath"
$ws1.Range("F1").Value = "iconPath"
$ws1.Range("G1").Value = "maxCount"
$ws1.Range("H1").Value = "isStat"
$ws1.Range("I1").Value = "isCrafting"

# --- Stats ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Stats")
$ws2.Range("A1").Value = "id"
$ws2.Range("B1").Value = "hp"
$ws2.Range("C1").Value = "temperature"
$ws2.Range("D1").Value = "atk"
$ws2.Range("E1").Value = "def"
$ws2.Range("F1").Value = "speed"
$ws2.Range("G1").Value = "stamina"

# --- Recipe ----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Recipe")
$ws3.Range("A1").Value = "id"
$ws3.Range("B1").Value = "craftingID"
$ws3.Range("C1").Value = "craftingPrice"
$ws3.Range("D1").Value = "materials_string"
$ws3.Range("E1").Value = "materials_count_string"
$ws3.Range("F1").Value = "availableCount"

# --- restore the saved selection / active-sheet UI state ------------------
$ws1.Range("I1").Select()
$ws2.Range("G1").Select()
$ws3.Range("E7").Select()
$ws3.Activate()
